# Apply the abstract edits described in the commit:
#  1. Add a clause about the magnetic field direction in the Introduction.
#  2. Mention the Gaussian wavepacket used with the finite difference method.
#  3. Add two blank paragraphs after the "Here the errors..." paragraph
#     (before "Literature cited").

$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# "...inhomogeneous (nonuniform) magnetic field the beam splits..."
#  -> "...magnetic field taken to be along the z-direction and with some
#      deviation in the x-direction the beam splits..."
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute(
    ") magnetic field the beam splits",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ") magnetic field taken to be along the z-direction and with some deviation in the x-direction the beam splits",
    2) | Out-Null

# --- Edit 2 -------------------------------------------------------------
# "...by using the finite difference method implemented in Python 3
#  software." -> "...by using the finite difference method with a
#  Gaussian wavepacket implemented in Python 3 software."
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute(
    "by using the finite difference method implemented in Python 3 software",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "by using the finite difference method with a Gaussian wavepacket implemented in Python 3 software",
    2) | Out-Null

# --- Edit 3 ---------------------------------------------------------------
# Insert two new empty paragraphs right after the paragraph that ends
# "...x-spatial derivatives." and before the existing blank paragraph that
# precedes "Literature cited".
$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Execute(
    "Here the errors are second ordered and the z",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hitPara = $find3.Parent.Paragraphs(1)
$hitPara.Range.InsertParagraphAfter() | Out-Null
$hitPara.Range.InsertParagraphAfter() | Out-Null
